# fix: capital letter in university name
# "СпбГУ" -> "СПбГУ" in the "название компании (если не СпбГУ)," cell,
# plus the two neighbouring label runs get re-applied (causing Word to
# coalesce the previously split "Консультант" + "," and "д" + "олжность,"
# runs), and the Normal style picks up an explicit suppressAutoHyphens
# toggle (same value already used for the document's paragraph defaults).

$d = $word.ActiveDocument

# Re-apply identical text across "Консультант" / "," so the two runs merge
# into a single run, matching how Word consolidates runs on a same-format
# edit.
$d.Content.Find.Execute("Консультант,", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Консультант,", 2)

# Same for "д" / "олжность,".
$d.Content.Find.Execute("должность,", $true, $false, $false, $false, $false, `
    $true, 1, $false, "должность,", 2)

# The actual content fix: lowercase "п" -> uppercase "П" in "СпбГУ".
$d.Content.Find.Execute("название компании (если не СпбГУ),", $true, $false, $false, $false, $false, `
    $true, 1, $false, "название компании (если не СПбГУ),", 2)

# Normal style: add an explicit suppressAutoHyphens toggle (already the
# value used by the document's paragraph defaults).
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Hyphenation = $false
